$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraphs we need by their text content (robust to
#    any index drift) instead of relying on fixed paragraph numbers.
# ------------------------------------------------------------------
$idxSeguidos = -1   # "Listagem de Pets seguidos por um Pet."
$idxAdicionar = -1  # "Adicionar evento."
$idxEventos = -1    # "Listagem de eventos."  (last match wins)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Listagem de Pets seguidos por um Pet*") {
        $idxSeguidos = $i
    }
    if ($t -like "*Adicionar evento*") {
        $idxAdicionar = $i
    }
    if ($t -like "*Listagem de eventos*") {
        $idxEventos = $i
    }
}

# The empty paragraph immediately following "Listagem de eventos." is
# also removed by the edit (it carried no numbering, just a stray
# formatted paragraph mark at the end of the "evento" block).
$idxBlankAfterEventos = $idxEventos + 1

# ------------------------------------------------------------------
# 2. Relocate the hidden "_GoBack" bookmark: it currently sits right
#    after "Comentar uma postagem." and needs to end up right after
#    "Listagem de Pets seguidos por um Pet." (before the trailing
#    space run we add below). Remove it from its old spot first.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3. Delete the four "evento" bullet paragraphs (Adicionar / Editar /
#    Deletar / Listagem de eventos) plus the trailing blank paragraph
#    right after them - this is the actual "remove the event" part of
#    the commit.
# ------------------------------------------------------------------
$pFrom = $d.Paragraphs.Item($idxAdicionar)
$pTo = $d.Paragraphs.Item($idxBlankAfterEventos)
$killRange = $d.Range($pFrom.Range.Start, $pTo.Range.End)
$killRange.Delete()

# ------------------------------------------------------------------
# 4. Append a trailing space run to the end of "Listagem de Pets
#    seguidos por um Pet." (the paragraph that now absorbs the old
#    stray formatting that used to trail the deleted "evento" block),
#    then re-create the "_GoBack" bookmark right before that new run.
# ------------------------------------------------------------------
$pSeguidos = $d.Paragraphs.Item($idxSeguidos)

# Insert point: right before this paragraph's own paragraph mark.
$insPoint = $d.Range($pSeguidos.Range.End - 1, $pSeguidos.Range.End - 1)
$insPoint.InsertAfter(" ")

# Turn the newly typed space into its own run (distinct w:r element)
# by toggling a character property on just that character.
$newRun = $d.Range($pSeguidos.Range.End - 2, $pSeguidos.Range.End - 1)
$newRun.Bold = 1
$newRun.Bold = 0

# Re-add the bookmark immediately before the new trailing space run.
$bmPoint = $d.Range($pSeguidos.Range.End - 2, $pSeguidos.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bmPoint)

Write-Output "done"
